$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84 - 李杰: 编写数据库设计文档 -> mark as 完成
$ws.Range("C84").Value = "完成"

# Row 85 - 周振朋: 编写用例规约 -> 编写"个人信息管理"用例规约, mark 完成了一部分
$ws.Range("B85").Value = "编写“个人信息管理”用例规约"
$ws.Range("C85").Value = "完成了一部分"
$ws.Rows.Item(85).RowHeight = 45

# Row 86 - 禤锦辉: 编写用例规约 -> 编写"首页"用例规约, mark 完成了一部分
$ws.Range("B86").Value = "编写“首页”用例规约"
$ws.Range("C86").Value = "完成了一部分"

# Row 87 - 柯新钿: 编写用例规约 -> 编写"账号管理"用例规约, mark 完成了一部分
$ws.Range("B87").Value = "编写“账号管理”用例规约"
$ws.Range("C87").Value = "完成了一部分"

# Row 88 - 冯文雄: 编写数据库设计文档 -> mark as 完成
$ws.Range("C88").Value = "完成"

# Row 89 - 阿卜力孜: 编写用例规约 -> 编写"查看买卖信息"用例规约, mark 完成了一部分
$ws.Range("B89").Value = "编写“查看买卖信息”用例规约"
$ws.Range("C89").Value = "完成了一部分"
$ws.Rows.Item(89).RowHeight = 45

# Row 90 - summary line
$ws.Range("A90").Value = "总结：明天早上马上完成用例规约"

# Move selection to reflect where the author left off editing
$ws.Range("B93").Select()
